$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# The sheet contains four identical "tables" (rows 4-8, 12-16, 20-26,
# 30-32) each with two placeholder columns (D:E) that only ever show an
# ellipsis "..." - a stand-in for a join column that was never filled
# in. Remove those two columns (shifting everything to their right one
# table-width to the left) across every table in one go.
$ws.Range("D4:E32").Delete(-4159)

# The delete-with-shift above also disturbs the small title block in
# row 2 (B2:D2), which sits just above the tables and is not supposed
# to be touched. Restore D2's original (blank, bold-row) formatting by
# copying it from its still-intact neighbour C2.
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the view state: scroll the window down so row 9 is at the top
# and select L29, matching where the author left off editing.
$ws.Range("L29").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
